$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary figures -----------------------------------------------
# "VALOR MORA" total (E11): 260000 -> 780000 (now covers 3 workers)
$ws.Range("E11").Value = 780000

# "Cant. Trabajadores" (C13): 1 -> 3 (two new workers added)
$ws.Range("C13").Value = 3

# --- Make room for the two new worker blocks (5 rows each) ----------------
# Insert 10 blank rows above the existing worker table (row 16), pushing the
# current CARLOS ALFREDO block (rows 16-20) and the footer (rows 25-26) down.
$ws.Range("A16:A25").EntireRow.Insert()

# The freshly inserted rows don't carry the table's row formatting yet -
# copy it from the row directly below (which still has the original style)
# down across the 10 new rows.
$ws.Range("B26:J26").Copy()
$ws.Range("B16:J25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Fill in the new worker: EDUI JOSE GARCIA VILLAMIZAR (rows 16-20) -----
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1052954524"
$ws.Range("D16").Value = "EDUI JOSE GARCIA VILLAMIZAR"
$ws.Range("E16").Value = "2502"
$ws.Range("F16").Value = 52000
$ws.Range("G16").Value = 1300000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1052954524"
$ws.Range("D17").Value = "EDUI JOSE GARCIA VILLAMIZAR"
$ws.Range("E17").Value = "2501"
$ws.Range("F17").Value = 52000
$ws.Range("G17").Value = 1300000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1052954524"
$ws.Range("D18").Value = "EDUI JOSE GARCIA VILLAMIZAR"
$ws.Range("E18").Value = "2412"
$ws.Range("F18").Value = 52000
$ws.Range("G18").Value = 1300000

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1052954524"
$ws.Range("D19").Value = "EDUI JOSE GARCIA VILLAMIZAR"
$ws.Range("E19").Value = "2411"
$ws.Range("F19").Value = 52000
$ws.Range("G19").Value = 1300000

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1052954524"
$ws.Range("D20").Value = "EDUI JOSE GARCIA VILLAMIZAR"
$ws.Range("E20").Value = "2410"
$ws.Range("F20").Value = 52000
$ws.Range("G20").Value = 1300000

# --- Fill in the new worker: MARIBEL BERMUDEZ DIAZ (rows 21-25) -----------
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "33065386"
$ws.Range("D21").Value = "MARIBEL BERMUDEZ DIAZ"
$ws.Range("E21").Value = "2502"
$ws.Range("F21").Value = 52000
$ws.Range("G21").Value = 689455

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "33065386"
$ws.Range("D22").Value = "MARIBEL BERMUDEZ DIAZ"
$ws.Range("E22").Value = "2501"
$ws.Range("F22").Value = 52000
$ws.Range("G22").Value = 689455

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "33065386"
$ws.Range("D23").Value = "MARIBEL BERMUDEZ DIAZ"
$ws.Range("E23").Value = "2412"
$ws.Range("F23").Value = 52000
$ws.Range("G23").Value = 689455

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "33065386"
$ws.Range("D24").Value = "MARIBEL BERMUDEZ DIAZ"
$ws.Range("E24").Value = "2411"
$ws.Range("F24").Value = 52000
$ws.Range("G24").Value = 689455

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "33065386"
$ws.Range("D25").Value = "MARIBEL BERMUDEZ DIAZ"
$ws.Range("E25").Value = "2410"
$ws.Range("F25").Value = 52000
$ws.Range("G25").Value = 689455

# --- Re-affirm the CARLOS ALFREDO RICARDO SOLA block (rows 26-30) ---------
# The row insert shifted this worker's original rows 16-20 down to 26-30,
# but the periods are re-listed in the same order used for the two new
# blocks above (2502, 2501, 2412, 2411, 2410), not the original ascending
# order, so set them explicitly rather than relying on the shifted values.
$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "1104380660"
$ws.Range("D26").Value = "CARLOS ALFREDO RICARDO SOLA"
$ws.Range("E26").Value = "2502"
$ws.Range("F26").Value = 52000
$ws.Range("G26").Value = 1300000

$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "1104380660"
$ws.Range("D27").Value = "CARLOS ALFREDO RICARDO SOLA"
$ws.Range("E27").Value = "2501"
$ws.Range("F27").Value = 52000
$ws.Range("G27").Value = 1300000

$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "1104380660"
$ws.Range("D28").Value = "CARLOS ALFREDO RICARDO SOLA"
$ws.Range("E28").Value = "2412"
$ws.Range("F28").Value = 52000
$ws.Range("G28").Value = 1300000

$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "1104380660"
$ws.Range("D29").Value = "CARLOS ALFREDO RICARDO SOLA"
$ws.Range("E29").Value = "2411"
$ws.Range("F29").Value = 52000
$ws.Range("G29").Value = 1300000

$ws.Range("B30").Value = "CC"
$ws.Range("C30").Value = "1104380660"
$ws.Range("D30").Value = "CARLOS ALFREDO RICARDO SOLA"
$ws.Range("E30").Value = "2410"
$ws.Range("F30").Value = 52000
$ws.Range("G30").Value = 1300000

# Rows 35-36 hold the footer signature block, auto-shifted down by the row
# insert above - nothing else to do there.
